# Add a new "2022-Q3" quarterly sheet (fund-holdings detail) right after the
# "总计" (totals) sheet, shifting 2022-Q2 / 2022-Q1 / 2021-Q3 one tab to the
# right, and roll the new quarter's aggregate numbers into "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet right after "总计" (becomes tab #2).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Populate "2022-Q3" with the fund-holding detail rows.
# ---------------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.BorderAround(1)
}

$rows = @(
    @("007340", "南方科技创新混合A",       "23.74", "79.07", "3.38", "0.8024", 10),
    @("007341", "南方科技创新混合C",       "4.54",  "79.07", "3.38", "0.1535", 10),
    @("001518", "万家瑞兴灵活配置混合A",   "2.81",  "83.43", "4.79", "0.1346", 3),
    @("020015", "国泰区位优势混合A",       "1.95",  "87.12", "5.40", "0.1053", 6),
    @("013869", "创金合信物联网主题股票A", "0.19",  "83.68", "8.12", "0.0154", 2),
    @("013870", "创金合信物联网主题股票C", "0.13",  "83.68", "8.12", "0.0106", 2),
    @("515510", "嘉实中证500成长估值ETF",   "0.13",  "98.46", "1.54", "0.0020", 1),
    @("515590", "前海开源中证500等权重ETF", "0.35",  "95.52", "0.33", "0.0012", 6),
    @("015594", "国泰区位优势混合C",       "0.00",  "87.12", "5.40", 0,        6),
    @("015390", "万家瑞兴灵活配置混合C",   "0.00",  "83.43", "4.79", 0,        3)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $idxCell = $q3.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.BorderAround(1)

    $q3.Cells.Item($excelRow, 2).Value = "'" + $row[0]
    $q3.Cells.Item($excelRow, 3).Value = "'" + $row[1]
    $q3.Cells.Item($excelRow, 4).Value = "'" + $row[2]
    $q3.Cells.Item($excelRow, 5).Value = "'" + $row[3]
    $q3.Cells.Item($excelRow, 6).Value = "'" + $row[4]

    $posValue = $row[5]
    if ($posValue -eq 0) {
        $q3.Cells.Item($excelRow, 7).Value = 0
    } else {
        $q3.Cells.Item($excelRow, 7).Value = "'" + $posValue
    }

    $q3.Cells.Item($excelRow, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 3. Roll the 2022-Q3 aggregate into the "总计" sheet: insert a fresh row
#    right below the header and push the older quarters down by one.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(5).Insert()
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(5, 2).Value = "2021-Q3"
$totalSheet.Cells.Item(5, 3).Value = 3
$totalSheet.Cells.Item(5, 4).Value = 0.24
$totalSheet.Cells.Item(5, 1).Font.Bold = $true
$totalSheet.Cells.Item(5, 1).HorizontalAlignment = -4108
$totalSheet.Cells.Item(5, 1).VerticalAlignment = -4160
$totalSheet.Cells.Item(5, 1).BorderAround(1)

$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 10
$totalSheet.Cells.Item(2, 4).Value = 1.23

$totalSheet.Cells.Item(3, 2).Value = "2022-Q2"
$totalSheet.Cells.Item(3, 3).Value = 2
$totalSheet.Cells.Item(3, 4).Value = 0.4

$totalSheet.Cells.Item(4, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(4, 3).Value = 2
$totalSheet.Cells.Item(4, 4).Value = 0

# ---------------------------------------------------------------------
# 4. Keep "2021-Q3" (now the last tab) the active/selected sheet, same as
#    before the edit.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
